# Update "想去人数" (F column) values for several events that changed between
# the previous site build and the one generated at commit 456a3b4.
#
# Sheet "展览":
#   F3  278  -> 279   (南宁·第五人格Only1.0)
#   F5  3406 -> 3408  (南宁·AP动漫游戏嘉年华)
#   F7  412  -> 413   (南宁·恋与深空only)
#   F10 35   -> 37    (广西·首届明日方舟only展 - 花庭圣梦)
#   F11 1261 -> 1262  (南宁·AB动漫游戏嘉年华)
#   F13 1534 -> 1543  (南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）)
#
# Sheet "全部类型" (same events, different row positions):
#   F3  278  -> 279
#   F5  3406 -> 3408
#   F7  412  -> 413
#   F11 35   -> 37
#   F14 1261 -> 1262
#   F16 1534 -> 1543

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 279
$wsExhibition.Range("F5").Value = 3408
$wsExhibition.Range("F7").Value = 413
$wsExhibition.Range("F10").Value = 37
$wsExhibition.Range("F11").Value = 1262
$wsExhibition.Range("F13").Value = 1543

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 279
$wsAll.Range("F5").Value = 3408
$wsAll.Range("F7").Value = 413
$wsAll.Range("F11").Value = 37
$wsAll.Range("F14").Value = 1262
$wsAll.Range("F16").Value = 1543
